$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5592.3667
$ws.Range("I132").Value = 2337.1155
$ws.Range("J132").Value = 26751.5
$ws.Range("K132").Value = 7011.3465
$ws.Range("L132").Value = 80254.5
$ws.Range("M132").Value = -4481.3465
$ws.Range("N132").Value = -85314.5
$ws.Range("H138").Value = 2328149.8
$ws.Range("I138").Value = 1235.7037
$ws.Range("J138").Value = 3393009
$ws.Range("K138").Value = 3707.1111
$ws.Range("L138").Value = 10179027
$ws.Range("M138").Value = 1432.8889
$ws.Range("N138").Value = -10189307
$ws.Range("H141").Value = 1520
$ws.Range("I141").Value = 930.5
$ws.Range("J141").Value = 3485
$ws.Range("K141").Value = 2791.5
$ws.Range("L141").Value = 10455
$ws.Range("M141").Value = 2388.5
$ws.Range("N141").Value = -20815

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2313.125
$ws.Range("I2").Value = 2328.182
$ws.Range("J2").Value = 2280
$ws.Range("K2").Value = 2328.182
$ws.Range("L2").Value = 2280
$ws.Range("M2").Value = -2215.182
$ws.Range("N2").Value = -2506
$ws.Range("H61").Value = 1343.4286
$ws.Range("I61").Value = 1205.3024
$ws.Range("J61").Value = 1800.3077
$ws.Range("K61").Value = 1205.3024
$ws.Range("L61").Value = 1800.3077
$ws.Range("M61").Value = -993.3024
$ws.Range("N61").Value = -2224.3077
$ws.Range("H74").Value = 22275.271
$ws.Range("I74").Value = 29657.629
$ws.Range("J74").Value = 2399.6924
$ws.Range("K74").Value = 29657.629
$ws.Range("L74").Value = 2399.6924
$ws.Range("M74").Value = -28783.629
$ws.Range("N74").Value = -4147.6924
$ws.Range("H77").Value = 22275.271
$ws.Range("I77").Value = 29657.629
$ws.Range("J77").Value = 2399.6924
$ws.Range("K77").Value = 148288.145
$ws.Range("L77").Value = 11998.462
$ws.Range("M77").Value = -143920.145
$ws.Range("N77").Value = -20734.462
$ws.Range("H113").Value = 30057.143
$ws.Range("J113").Value = 30057.143
$ws.Range("L113").Value = 30057.143
$ws.Range("N113").Value = -38735.143
$ws.Range("H116").Value = 2313.125
$ws.Range("I116").Value = 2328.182
$ws.Range("J116").Value = 2280
$ws.Range("K116").Value = 2328.182
$ws.Range("L116").Value = 2280
$ws.Range("M116").Value = -34.18199999999979
$ws.Range("N116").Value = -6868
$ws.Range("H132").Value = 1415.7742
$ws.Range("I132").Value = 1166.85
$ws.Range("J132").Value = 1868.3636
$ws.Range("K132").Value = 3500.55
$ws.Range("L132").Value = 5605.0908
$ws.Range("M132").Value = -970.5499999999997
$ws.Range("N132").Value = -10665.0908
$ws.Range("H136").Value = 1343.4286
$ws.Range("I136").Value = 1205.3024
$ws.Range("J136").Value = 1800.3077
$ws.Range("K136").Value = 3615.9072
$ws.Range("L136").Value = 5400.9231
$ws.Range("M136").Value = -1065.9072
$ws.Range("N136").Value = -10500.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2313.125
$ws.Range("I3").Value = 2328.182
$ws.Range("J3").Value = 2280
$ws.Range("K3").Value = 2328.182
$ws.Range("L3").Value = 2280
$ws.Range("M3").Value = -2214.182
$ws.Range("N3").Value = -2508
$ws.Range("H134").Value = 609516.9
$ws.Range("I134").Value = 1146022
$ws.Range("J134").Value = 3785.1936
$ws.Range("K134").Value = 3438066
$ws.Range("L134").Value = 11355.5808
$ws.Range("M134").Value = -3435531
$ws.Range("N134").Value = -16425.5808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1961.0492
$ws.Range("I31").Value = 1221.6666
$ws.Range("J31").Value = 2440.6487
$ws.Range("K31").Value = 1221.6666
$ws.Range("L31").Value = 2440.6487
$ws.Range("M31").Value = -926.6666
$ws.Range("N31").Value = -3030.6487
$ws.Range("H34").Value = 1961.0492
$ws.Range("I34").Value = 1221.6666
$ws.Range("J34").Value = 2440.6487
$ws.Range("K34").Value = 1221.6666
$ws.Range("L34").Value = 2440.6487
$ws.Range("M34").Value = -1019.6666
$ws.Range("N34").Value = -2844.6487
$ws.Range("H58").Value = 3009.0667
$ws.Range("I58").Value = 3655.5881
$ws.Range("J58").Value = 1010.7273
$ws.Range("K58").Value = 3655.5881
$ws.Range("L58").Value = 1010.7273
$ws.Range("M58").Value = -3452.5881
$ws.Range("N58").Value = -1416.7273
$ws.Range("H132").Value = 640030.5600000001
$ws.Range("I132").Value = 1342
$ws.Range("J132").Value = 3705735.8
$ws.Range("K132").Value = 4026
$ws.Range("L132").Value = 11117207.4
$ws.Range("M132").Value = -1496
$ws.Range("N132").Value = -11122267.4
$ws.Range("H134").Value = 1527.4875
$ws.Range("I134").Value = 1495.661
$ws.Range("J134").Value = 1616.9048
$ws.Range("K134").Value = 4486.983
$ws.Range("L134").Value = 4850.7144
$ws.Range("M134").Value = -1951.983
$ws.Range("N134").Value = -9920.714400000001
$ws.Range("H136").Value = 3009.0667
$ws.Range("I136").Value = 3655.5881
$ws.Range("J136").Value = 1010.7273
$ws.Range("K136").Value = 10966.7643
$ws.Range("L136").Value = 3032.1819
$ws.Range("M136").Value = -8416.764299999999
$ws.Range("N136").Value = -8132.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 888
$ws.Range("I97").Value = 417.64706
$ws.Range("J97").Value = 1221.1666
$ws.Range("K97").Value = 1252.94118
$ws.Range("L97").Value = 3663.4998
$ws.Range("M97").Value = -756.94118
$ws.Range("N97").Value = -4655.4998
$ws.Range("H113").Value = 1443426.8
$ws.Range("I113").Value = 3788361.5
$ws.Range("J113").Value = 389.92307
$ws.Range("K113").Value = 11365084.5
$ws.Range("L113").Value = 1169.76921
$ws.Range("M113").Value = -11362914.5
$ws.Range("N113").Value = -5509.76921
$ws.Range("H131").Value = 918.27
$ws.Range("J131").Value = 921.7217000000001
$ws.Range("L131").Value = 2765.1651
$ws.Range("N131").Value = -12845.1651

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1697428.5
$ws.Range("I132").Value = 2276.5527
$ws.Range("J132").Value = 4764846
$ws.Range("K132").Value = 6829.658100000001
$ws.Range("L132").Value = 14294538
$ws.Range("M132").Value = -4299.658100000001
$ws.Range("N132").Value = -14299598

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1766.1082
$ws.Range("I40").Value = 1655.5
$ws.Range("J40").Value = 2027.5454
$ws.Range("K40").Value = 1655.5
$ws.Range("L40").Value = 2027.5454
$ws.Range("M40").Value = -1519.5
$ws.Range("N40").Value = -2299.5454
$ws.Range("H88").Value = 29000
$ws.Range("I88").Value = 29000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 29000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -28572
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 29000
$ws.Range("I91").Value = 29000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 29000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -27518
$ws.Range("N91").ClearContents()
$ws.Range("H132").Value = 3615.8076
$ws.Range("I132").Value = 3407.027
$ws.Range("J132").Value = 4130.8
$ws.Range("K132").Value = 10221.081
$ws.Range("L132").Value = 12392.4
$ws.Range("M132").Value = -7691.081
$ws.Range("N132").Value = -17452.4
$ws.Range("H136").Value = 1385.9822
$ws.Range("I136").Value = 1025.5555
$ws.Range("J136").Value = 2034.75
$ws.Range("K136").Value = 3076.6665
$ws.Range("L136").Value = 6104.25
$ws.Range("M136").Value = -526.6664999999998
$ws.Range("N136").Value = -11204.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2397.7737
$ws.Range("I136").Value = 2354.2122
$ws.Range("J136").Value = 2469.65
$ws.Range("K136").Value = 7062.6366
$ws.Range("L136").Value = 7408.950000000001
$ws.Range("M136").Value = -4512.6366
$ws.Range("N136").Value = -12508.95
